# Refresh the crypto Price (column D) and Volume(1h) (column E) figures
# in the cryptos worksheet with the latest data.
#
# Some Price values are plain decimal numbers (e.g. "518.60"); setting
# .Value directly on those would let Excel auto-coerce them into numeric
# cells and silently drop the trailing zero (518.60 -> 518.6). To keep
# them as text exactly as scraped, the cell is temporarily switched to a
# text number format before the assignment, then its format/style is put
# back to the default afterward.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "56.450.59"
$ws.Range("E2").Value = "  +2.15%  "
$ws.Range("D3").Value = "2.315.76"
$ws.Range("E3").Value = "  +0.86%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "518.60"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.21%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "134.27"
$ws.Range("D6").NumberFormat = "General"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.54%  "
$ws.Range("E7").Value = "  +0.20%  "
$ws.Range("E8").Value = "  +0.92%  "
$ws.Range("D9").Value = "2.337.13"
$ws.Range("E9").Value = "  +0.80%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.103"
$ws.Range("D10").NumberFormat = "General"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.46%  "
$ws.Range("E11").Value = "  -1.06%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.30"
$ws.Range("D12").NumberFormat = "General"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +4.09%  "
$ws.Range("E13").Value = "  +0.05%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "23.85"
$ws.Range("D14").NumberFormat = "General"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.59%  "
$ws.Range("D15").Value = "2.730.27"
$ws.Range("E15").Value = "  +1.08%  "
$ws.Range("D16").Value = "56.568.95"
$ws.Range("E16").Value = "  +2.32%  "
$ws.Range("E17").Value = "  +1.59%  "
$ws.Range("D18").Value = "2.334.11"
$ws.Range("E18").Value = "  +1.35%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.45"
$ws.Range("D19").NumberFormat = "General"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.20%  "
$ws.Range("E20").Value = "  +0.95%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "323.18"
$ws.Range("D21").NumberFormat = "General"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.32%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.55"
$ws.Range("D22").NumberFormat = "General"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.56%  "
$ws.Range("E23").Value = "  +0.14%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "60.78"
$ws.Range("D24").NumberFormat = "General"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.85%  "
$ws.Range("E25").Value = "  +5.22%  "
$ws.Range("E26").Value = "  +0.02%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.96"
$ws.Range("D27").NumberFormat = "General"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +5.78%  "
$ws.Range("E28").Value = "  +10.95%  "
$ws.Range("D29").Value = "0.0₃0738"
$ws.Range("E29").Value = "  +3.50%  "
$ws.Range("E30").Value = "  +3.55%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "166.33"
$ws.Range("D31").NumberFormat = "General"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.52%  "
$ws.Range("E32").Value = "  +0.21%  "
$ws.Range("E33").Value = "  +1.66%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.992"
$ws.Range("D35").NumberFormat = "General"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.00%  "
$ws.Range("E36").Value = "  +0.96%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.921"
$ws.Range("D37").NumberFormat = "General"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.11%  "
$ws.Range("E38").Value = "  +3.07%  "
$ws.Range("E39").Value = "  +5.64%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "37.92"
$ws.Range("D40").NumberFormat = "General"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.78%  "
$ws.Range("E41").Value = "  +1.19%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "139.18"
$ws.Range("D42").NumberFormat = "General"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.55%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.60"
$ws.Range("D43").NumberFormat = "General"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.91%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.24"
$ws.Range("D44").NumberFormat = "General"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.60%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "278.85"
$ws.Range("D45").NumberFormat = "General"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +6.27%  "
$ws.Range("E46").Value = "  +1.60%  "
$ws.Range("E47").Value = "  -0.37%  "
$ws.Range("E48").Value = "  +1.49%  "
$ws.Range("E49").Value = "  +2.40%  "
$ws.Range("E50").Value = "  +1.01%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "17.82"
$ws.Range("D51").NumberFormat = "General"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +7.80%  "
